$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Snr" (column J header) is being dropped from the table. I1 ("Dispersion")
# becomes the new right-most bordered header cell, so give it J1's current
# header formatting (fill + border + alignment) before clearing J1 out.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null

# J1 becomes an empty filler cell: no text, no border, just a light fill.
$ws.Range("J1").Value = ""
$ws.Range("J1").Borders.LineStyle = -4142
$ws.Range("J1").Interior.ThemeColor = 2

# Move the active selection to I5.
$ws.Range("I5").Select() | Out-Null
